$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Scratch cell used to stamp each new tracking number in as TEXT (so the
# written cell keeps a shared-string type, matching the original cell,
# without ever changing the destination cell's own style/number format).
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"

function Set-TrackingNumber($cellRef, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163) # xlPasteValues
}

Set-TrackingNumber "C2" "320018799382"
Set-TrackingNumber "C3" "320018799393"
Set-TrackingNumber "C4" "320018799420"
Set-TrackingNumber "C5" "320018799441"
Set-TrackingNumber "D5" "320018799441"
Set-TrackingNumber "C6" "320018799485"
Set-TrackingNumber "D6" "320018799485"
Set-TrackingNumber "C7" "320018799500"
Set-TrackingNumber "D7" "320018799500"
Set-TrackingNumber "C8" "320018799533"
Set-TrackingNumber "C9" "320018799625"
Set-TrackingNumber "C10" "320018799658"
Set-TrackingNumber "C11" "320018799670"
Set-TrackingNumber "C12" "320018799717"
Set-TrackingNumber "C13" "320018799739"
Set-TrackingNumber "D13" "320018799739"
Set-TrackingNumber "C14" "320018799761"
Set-TrackingNumber "D14" "320018799761"
Set-TrackingNumber "C15" "320018799783"
Set-TrackingNumber "D15" "320018799783"
Set-TrackingNumber "C16" "320018799810"
Set-TrackingNumber "D16" "320018799810"
Set-TrackingNumber "C17" "320018799831"
Set-TrackingNumber "D17" "320018799831"
Set-TrackingNumber "C18" "320018799875"
Set-TrackingNumber "C19" "320018792701"
Set-TrackingNumber "C20" "320018792734"
Set-TrackingNumber "C21" "320018792756"
Set-TrackingNumber "C22" "320018792789"

$scratch.Clear()
$excel.CutCopyMode = 0

